{"js": "// Apply the text replacements described by the diff: update the date\n// line and each \"a\u00f7b=\" division expression in the practice table.\nconst replacements = [\n  [\"2024-01-20 Saturday\", \"2024-01-21 Sunday\"],\n  [\"315\u00f75=\", \"345\u00f73=\"],\n  [\"342\u00f77=\", \"160\u00f74=\"],\n  [\"860\u00f78=\", \"152\u00f73=\"],\n  [\"622\u00f75=\", \"654\u00f76=\"],\n  [\"865\u00f72=\", \"257\u00f72=\"],\n  [\"902\u00f79=\", \"215\u00f75=\"],\n  [\"302\u00f78=\", \"184\u00f79=\"],\n  [\"614\u00f79=\", \"828\u00f78=\"],\n  [\"131\u00f75=\", \"601\u00f77=\"],\n  [\"739\u00f72=\", \"900\u00f75=\"],\n  [\"869\u00f77=\", \"835\u00f78=\"],\n  [\"723\u00f74=\", \"713\u00f79=\"],\n  [\"557\u00f76=\", \"551\u00f79=\"],\n  [\"325\u00f79=\", \"422\u00f79=\"],\n  [\"536\u00f75=\", \"778\u00f72=\"],\n  [\"811\u00f73=\", \"398\u00f78=\"],\n  [\"326\u00f75=\", \"917\u00f77=\"],\n  [\"376\u00f77=\", \"527\u00f72=\"],\n  [\"918\u00f75=\", \"471\u00f75=\"],\n  [\"107\u00f76=\", \"448\u00f74=\"],\n  [\"698\u00f73=\", \"744\u00f76=\"],\n  [\"231\u00f79=\", \"414\u00f73=\"],\n  [\"443\u00f74=\", \"405\u00f74=\"],\n  [\"805\u00f78=\", \"940\u00f74=\"],\n  [\"223\u00f74=\", \"116\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff: update the date\n# line and each \"a\u00f7b=\" division expression in the practice table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-01-20 Saturday\"; New = \"2024-01-21 Sunday\" },\n    @{ Old = \"315\u00f75=\"; New = \"345\u00f73=\" },\n    @{ Old = \"342\u00f77=\"; New = \"160\u00f74=\" },\n    @{ Old = \"860\u00f78=\"; New = \"152\u00f73=\" },\n    @{ Old = \"622\u00f75=\"; New = \"654\u00f76=\" },\n    @{ Old = \"865\u00f72=\"; New = \"257\u00f72=\" },\n    @{ Old = \"902\u00f79=\"; New = \"215\u00f75=\" },\n    @{ Old = \"302\u00f78=\"; New = \"184\u00f79=\" },\n    @{ Old = \"614\u00f79=\"; New = \"828\u00f78=\" },\n    @{ Old = \"131\u00f75=\"; New = \"601\u00f77=\" },\n    @{ Old = \"739\u00f72=\"; New = \"900\u00f75=\" },\n    @{ Old = \"869\u00f77=\"; New = \"835\u00f78=\" },\n    @{ Old = \"723\u00f74=\"; New = \"713\u00f79=\" },\n    @{ Old = \"557\u00f76=\"; New = \"551\u00f79=\" },\n    @{ Old = \"325\u00f79=\"; New = \"422\u00f79=\" },\n    @{ Old = \"536\u00f75=\"; New = \"778\u00f72=\" },\n    @{ Old = \"811\u00f73=\"; New = \"398\u00f78=\" },\n    @{ Old = \"326\u00f75=\"; New = \"917\u00f77=\" },\n    @{ Old = \"376\u00f77=\"; New = \"527\u00f72=\" },\n    @{ Old = \"918\u00f75=\"; New = \"471\u00f75=\" },\n    @{ Old = \"107\u00f76=\"; New = \"448\u00f74=\" },\n    @{ Old = \"698\u00f73=\"; New = \"744\u00f76=\" },\n    @{ Old = \"231\u00f79=\"; New = \"414\u00f73=\" },\n    @{ Old = \"443\u00f74=\"; New = \"405\u00f74=\" },\n    @{ Old = \"805\u00f78=\"; New = \"940\u00f74=\" },\n    @{ Old = \"223\u00f74=\"; New = \"116\u00f74=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
